$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntrk3"
$ws.Range("C2").Value = "Ptprf"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1.0
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.004481
$ws.Range("H2").Value = 0.013443
$ws.Range("I2").Value = 0.01111165297720057
$ws.Range("J2").Value = 0.01111165297720057
$ws.Range("K2").Value = 3.0
$ws.Range("L2").Value = 1.0
$ws.Range("M2").Value = 0.1863256666666667
$ws.Range("N2").Value = 0.5589770000000001
$ws.Range("O2").Value = 0.01657678358851065
$ws.Range("P2").Value = 0.01657678358851065
$ws.Range("Q2").Value = 0.0008349253123333334
$ws.Range("R2").Value = 0.007514327811000001
$ws.Range("S2").Value = 0.0001841954667136839
$ws.Range("T2").Value = 0.0001841954667136839

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ntrk3"
$ws.Range("C3").Value = "Ptprf"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1.0
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.004481
$ws.Range("H3").Value = 0.013443
$ws.Range("I3").Value = 0.01111165297720057
$ws.Range("J3").Value = 0.01111165297720057
$ws.Range("K3").Value = 3.0
$ws.Range("L3").Value = 1.0
$ws.Range("M3").Value = 5.829902333333333
$ws.Range("N3").Value = 17.489707
$ws.Range("O3").Value = 0.5186672939413604
$ws.Range("P3").Value = 0.5186672939413604
$ws.Range("Q3").Value = 0.02612379235566666
$ws.Range("R3").Value = 0.235114131201
$ws.Range("S3").Value = 0.00576325098090008
$ws.Range("T3").Value = 0.005763250980900081

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ntrk3"
$ws.Range("C4").Value = "Ptprf"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1.0
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.004481
$ws.Range("H4").Value = 0.013443
$ws.Range("I4").Value = 0.01111165297720057
$ws.Range("J4").Value = 0.01111165297720057
$ws.Range("K4").Value = 3.0
$ws.Range("L4").Value = 1.0
$ws.Range("M4").Value = 5.212463666666667
$ws.Range("N4").Value = 15.637391
$ws.Range("O4").Value = 0.4637358003923671
$ws.Range("P4").Value = 0.4637358003923669
$ws.Range("Q4").Value = 0.02335704969033333
$ws.Range("R4").Value = 0.210213447213
$ws.Range("S4").Value = 0.005152871287064335
$ws.Range("T4").Value = 0.005152871287064334

# Row 5: ECs -> Resolving-Mac
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Ntrk3"
$ws.Range("C5").Value = "Ptprf"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1.0
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.004481
$ws.Range("H5").Value = 0.013443
$ws.Range("I5").Value = 0.01111165297720057
$ws.Range("J5").Value = 0.01111165297720057
$ws.Range("K5").Value = 3.0
$ws.Range("L5").Value = 1.0
$ws.Range("M5").Value = 0.01146633333333333
$ws.Range("N5").Value = 0.034399
$ws.Range("O5").Value = 0.001020122077762015
$ws.Range("P5").Value = 0.001020122077762015
$ws.Range("Q5").Value = 0.00005138063966666667
$ws.Range("R5").Value = 0.000462425757
$ws.Range("S5").Value = 0.00001133524252247233
$ws.Range("T5").Value = 0.00001133524252247233

# Row 6: FAPs -> ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ntrk3"
$ws.Range("C6").Value = "Ptprf"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3.0
$ws.Range("F6").Value = 1.0
$ws.Range("G6").Value = 0.1511493333333333
$ws.Range("H6").Value = 0.453448
$ws.Range("I6").Value = 0.3748089577628241
$ws.Range("J6").Value = 0.3748089577628241
$ws.Range("K6").Value = 3.0
$ws.Range("L6").Value = 1.0
$ws.Range("M6").Value = 0.1863256666666667
$ws.Range("N6").Value = 0.5589770000000001
$ws.Range("O6").Value = 0.01657678358851065
$ws.Range("P6").Value = 0.01657678358851065
$ws.Range("Q6").Value = 0.02816300029955556
$ws.Range("R6").Value = 0.253467002696
$ws.Range("S6").Value = 0.006213126979869564
$ws.Range("T6").Value = 0.006213126979869562

# Row 7: FAPs -> FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ntrk3"
$ws.Range("C7").Value = "Ptprf"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3.0
$ws.Range("F7").Value = 1.0
$ws.Range("G7").Value = 0.1511493333333333
$ws.Range("H7").Value = 0.453448
$ws.Range("I7").Value = 0.3748089577628241
$ws.Range("J7").Value = 0.3748089577628241
$ws.Range("K7").Value = 3.0
$ws.Range("L7").Value = 1.0
$ws.Range("M7").Value = 5.829902333333333
$ws.Range("N7").Value = 17.489707
$ws.Range("O7").Value = 0.5186672939413604
$ws.Range("P7").Value = 0.5186672939413604
$ws.Range("Q7").Value = 0.8811858510817777
$ws.Range("R7").Value = 7.930672659735999
$ws.Range("S7").Value = 0.1944011478678256
$ws.Range("T7").Value = 0.1944011478678256

# Row 8: FAPs -> MuSCs
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Ntrk3"
$ws.Range("C8").Value = "Ptprf"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3.0
$ws.Range("F8").Value = 1.0
$ws.Range("G8").Value = 0.1511493333333333
$ws.Range("H8").Value = 0.453448
$ws.Range("I8").Value = 0.3748089577628241
$ws.Range("J8").Value = 0.3748089577628241
$ws.Range("K8").Value = 3.0
$ws.Range("L8").Value = 1.0
$ws.Range("M8").Value = 5.212463666666667
$ws.Range("N8").Value = 15.637391
$ws.Range("O8").Value = 0.4637358003923671
$ws.Range("P8").Value = 0.4637358003923669
$ws.Range("Q8").Value = 0.787860408240889
$ws.Range("R8").Value = 7.090743674167999
$ws.Range("S8").Value = 0.1738123320223721
$ws.Range("T8").Value = 0.1738123320223721

# Row 9: FAPs -> Resolving-Mac
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Ntrk3"
$ws.Range("C9").Value = "Ptprf"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3.0
$ws.Range("F9").Value = 1.0
$ws.Range("G9").Value = 0.1511493333333333
$ws.Range("H9").Value = 0.453448
$ws.Range("I9").Value = 0.3748089577628241
$ws.Range("J9").Value = 0.3748089577628241
$ws.Range("K9").Value = 3.0
$ws.Range("L9").Value = 1.0
$ws.Range("M9").Value = 0.01146633333333333
$ws.Range("N9").Value = 0.034399
$ws.Range("O9").Value = 0.001020122077762015
$ws.Range("P9").Value = 0.001020122077762015
$ws.Range("Q9").Value = 0.001733128639111111
$ws.Range("R9").Value = 0.015598157752
$ws.Range("S9").Value = 0.0003823508927568274
$ws.Range("T9").Value = 0.0003823508927568273

# Row 10: MuSCs -> ECs
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Ntrk3"
$ws.Range("C10").Value = "Ptprf"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3.0
$ws.Range("F10").Value = 1.0
$ws.Range("G10").Value = 0.24764
$ws.Range("H10").Value = 0.74292
$ws.Range("I10").Value = 0.6140793892599753
$ws.Range("J10").Value = 0.6140793892599754
$ws.Range("K10").Value = 3.0
$ws.Range("L10").Value = 1.0
$ws.Range("M10").Value = 0.1863256666666667
$ws.Range("N10").Value = 0.5589770000000001
$ws.Range("O10").Value = 0.01657678358851065
$ws.Range("P10").Value = 0.01657678358851065
$ws.Range("Q10").Value = 0.04614168809333334
$ws.Range("R10").Value = 0.41527519284
$ws.Range("S10").Value = 0.0101794611419274
$ws.Range("T10").Value = 0.0101794611419274

# Row 11: MuSCs -> FAPs
$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Ntrk3"
$ws.Range("C11").Value = "Ptprf"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3.0
$ws.Range("F11").Value = 1.0
$ws.Range("G11").Value = 0.24764
$ws.Range("H11").Value = 0.74292
$ws.Range("I11").Value = 0.6140793892599753
$ws.Range("J11").Value = 0.6140793892599754
$ws.Range("K11").Value = 3.0
$ws.Range("L11").Value = 1.0
$ws.Range("M11").Value = 5.829902333333333
$ws.Range("N11").Value = 17.489707
$ws.Range("O11").Value = 0.5186672939413604
$ws.Range("P11").Value = 0.5186672939413604
$ws.Range("Q11").Value = 1.443717013826667
$ws.Range("R11").Value = 12.99345312444
$ws.Range("S11").Value = 0.3185028950926346
$ws.Range("T11").Value = 0.3185028950926347

# Row 12: MuSCs -> MuSCs
$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Ntrk3"
$ws.Range("C12").Value = "Ptprf"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3.0
$ws.Range("F12").Value = 1.0
$ws.Range("G12").Value = 0.24764
$ws.Range("H12").Value = 0.74292
$ws.Range("I12").Value = 0.6140793892599753
$ws.Range("J12").Value = 0.6140793892599754
$ws.Range("K12").Value = 3.0
$ws.Range("L12").Value = 1.0
$ws.Range("M12").Value = 5.212463666666667
$ws.Range("N12").Value = 15.637391
$ws.Range("O12").Value = 0.4637358003923671
$ws.Range("P12").Value = 0.4637358003923669
$ws.Range("Q12").Value = 1.290814502413334
$ws.Range("R12").Value = 11.61733052172
$ws.Range("S12").Value = 0.2847705970829306
$ws.Range("T12").Value = 0.2847705970829306

# Row 13: MuSCs -> Resolving-Mac
$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Ntrk3"
$ws.Range("C13").Value = "Ptprf"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3.0
$ws.Range("F13").Value = 1.0
$ws.Range("G13").Value = 0.24764
$ws.Range("H13").Value = 0.74292
$ws.Range("I13").Value = 0.6140793892599753
$ws.Range("J13").Value = 0.6140793892599754
$ws.Range("K13").Value = 3.0
$ws.Range("L13").Value = 1.0
$ws.Range("M13").Value = 0.01146633333333333
$ws.Range("N13").Value = 0.034399
$ws.Range("O13").Value = 0.001020122077762015
$ws.Range("P13").Value = 0.001020122077762015
$ws.Range("Q13").Value = 0.002839522786666667
$ws.Range("R13").Value = 0.02555570508
$ws.Range("S13").Value = 0.0006264359424827152
$ws.Range("T13").Value = 0.0006264359424827152
